# Auto-generated edit script: updates Leve market-price columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled
# market-data refresh run.

$wb = $excel.ActiveWorkbook

$edits = @(
    ,@("ALC", "H3", "set", 35000)
    ,@("ALC", "J3", "set", 35000)
    ,@("ALC", "L3", "set", 35000)
    ,@("ALC", "N3", "set", -35228)
    ,@("ALC", "H7", "set", 10005)
    ,@("ALC", "J7", "set", 0)
    ,@("ALC", "L7", "set", 0)
    ,@("ALC", "N7", "remove", $null)
    ,@("ALC", "H14", "set", 10005)
    ,@("ALC", "J14", "set", 0)
    ,@("ALC", "L14", "set", 0)
    ,@("ALC", "N14", "remove", $null)
    ,@("ALC", "H96", "set", 270.53333)
    ,@("ALC", "I96", "set", 240.6923)
    ,@("ALC", "J96", "set", 464.5)
    ,@("ALC", "K96", "set", 722.0769)
    ,@("ALC", "L96", "set", 1393.5)
    ,@("ALC", "M96", "set", 650.9231)
    ,@("ALC", "N96", "set", -4139.5)
    ,@("ALC", "H98", "set", 1018)
    ,@("ALC", "I98", "set", 1040.5385)
    ,@("ALC", "J98", "set", 725)
    ,@("ALC", "K98", "set", 1040.5385)
    ,@("ALC", "L98", "set", 725)
    ,@("ALC", "M98", "set", 457.4614999999999)
    ,@("ALC", "N98", "set", -3721)
    ,@("ALC", "H102", "set", 35000)
    ,@("ALC", "J102", "set", 35000)
    ,@("ALC", "L102", "set", 35000)
    ,@("ALC", "N102", "set", -41490)
    ,@("ALC", "H111", "set", 1176.375)
    ,@("ALC", "I111", "set", 1063.1666)
    ,@("ALC", "J111", "set", 1516)
    ,@("ALC", "K111", "set", 3189.4998)
    ,@("ALC", "L111", "set", 4548)
    ,@("ALC", "M111", "set", -122.4998000000001)
    ,@("ALC", "N111", "set", -10682)
    ,@("ALC", "H122", "set", 1018)
    ,@("ALC", "I122", "set", 1040.5385)
    ,@("ALC", "J122", "set", 725)
    ,@("ALC", "K122", "set", 3121.6155)
    ,@("ALC", "L122", "set", 2175)
    ,@("ALC", "M122", "set", -671.6155000000003)
    ,@("ALC", "N122", "set", -7075)
    ,@("ALC", "H127", "set", 1582.75)
    ,@("ALC", "J127", "set", 2037)
    ,@("ALC", "L127", "set", 6111)
    ,@("ALC", "N127", "set", -16031)
    ,@("ALC", "H137", "set", 2859709.8)
    ,@("ALC", "I137", "set", 4763698.5)
    ,@("ALC", "J137", "set", 3726.7856)
    ,@("ALC", "K137", "set", 14291095.5)
    ,@("ALC", "L137", "set", 11180.3568)
    ,@("ALC", "M137", "set", -14288545.5)
    ,@("ALC", "N137", "set", -16280.3568)
    ,@("ALC", "H138", "set", 3087936)
    ,@("ALC", "I138", "set", 964.0270400000001)
    ,@("ALC", "J138", "set", 9806640)
    ,@("ALC", "K138", "set", 2892.08112)
    ,@("ALC", "L138", "set", 29419920)
    ,@("ALC", "M138", "set", 2247.91888)
    ,@("ALC", "N138", "set", -29430200)
    ,@("ARM", "H31", "set", 0)
    ,@("ARM", "I31", "set", 0)
    ,@("ARM", "K31", "set", 0)
    ,@("ARM", "M31", "remove", $null)
    ,@("ARM", "H69", "set", 47229.5)
    ,@("ARM", "J69", "set", 47229.5)
    ,@("ARM", "L69", "set", 47229.5)
    ,@("ARM", "N69", "set", -48727.5)
    ,@("ARM", "H72", "set", 47229.5)
    ,@("ARM", "J72", "set", 47229.5)
    ,@("ARM", "L72", "set", 141688.5)
    ,@("ARM", "N72", "set", -149176.5)
    ,@("ARM", "H93", "set", 13000)
    ,@("ARM", "J93", "set", 13000)
    ,@("ARM", "L93", "set", 13000)
    ,@("ARM", "N93", "set", -17992)
    ,@("CRP", "H17", "set", 0)
    ,@("CRP", "I17", "set", 0)
    ,@("CRP", "K17", "set", 0)
    ,@("CRP", "M17", "remove", $null)
    ,@("CRP", "H31", "set", 1623.3247)
    ,@("CRP", "I31", "set", 1107.0847)
    ,@("CRP", "J31", "set", 3315.4443)
    ,@("CRP", "K31", "set", 1107.0847)
    ,@("CRP", "L31", "set", 3315.4443)
    ,@("CRP", "M31", "set", -812.0847000000001)
    ,@("CRP", "N31", "set", -3905.4443)
    ,@("CRP", "H34", "set", 1623.3247)
    ,@("CRP", "I34", "set", 1107.0847)
    ,@("CRP", "J34", "set", 3315.4443)
    ,@("CRP", "K34", "set", 1107.0847)
    ,@("CRP", "L34", "set", 3315.4443)
    ,@("CRP", "M34", "set", -905.0847000000001)
    ,@("CRP", "N34", "set", -3719.4443)
    ,@("CRP", "H122", "set", 2378.1667)
    ,@("CRP", "I122", "set", 1804.1333)
    ,@("CRP", "K122", "set", 5412.3999)
    ,@("CRP", "M122", "set", -2962.3999)
    ,@("CRP", "H134", "set", 29119.975)
    ,@("CRP", "I134", "set", 1560.6129)
    ,@("CRP", "J134", "set", 135912.5)
    ,@("CRP", "K134", "set", 4681.8387)
    ,@("CRP", "L134", "set", 407737.5)
    ,@("CRP", "M134", "set", -2146.8387)
    ,@("CRP", "N134", "set", -412807.5)
    ,@("CUL", "H4", "set", 3228169.8)
    ,@("CUL", "J4", "set", 3335712)
    ,@("CUL", "L4", "set", 10007136)
    ,@("CUL", "N4", "set", -10007360)
    ,@("CUL", "H23", "set", 141)
    ,@("CUL", "I23", "set", 95)
    ,@("CUL", "J23", "set", 164)
    ,@("CUL", "K23", "set", 285)
    ,@("CUL", "L23", "set", 492)
    ,@("CUL", "M23", "set", -50)
    ,@("CUL", "N23", "set", -962)
    ,@("GSM", "H3", "set", 1901)
    ,@("GSM", "I3", "set", 1851.5)
    ,@("GSM", "K3", "set", 1851.5)
    ,@("GSM", "M3", "set", -1735.5)
    ,@("GSM", "H34", "set", 35173)
    ,@("GSM", "J34", "set", 35173)
    ,@("GSM", "L34", "set", 35173)
    ,@("GSM", "N34", "set", -35709)
    ,@("GSM", "H70", "set", 5358.8125)
    ,@("GSM", "I70", "set", 5137.2856)
    ,@("GSM", "J70", "set", 5781.727)
    ,@("GSM", "K70", "set", 5137.2856)
    ,@("GSM", "L70", "set", 5781.727)
    ,@("GSM", "M70", "set", -4867.2856)
    ,@("GSM", "N70", "set", -6321.727)
    ,@("GSM", "H73", "set", 5358.8125)
    ,@("GSM", "I73", "set", 5137.2856)
    ,@("GSM", "J73", "set", 5781.727)
    ,@("GSM", "K73", "set", 5137.2856)
    ,@("GSM", "L73", "set", 5781.727)
    ,@("GSM", "M73", "set", -4201.2856)
    ,@("GSM", "N73", "set", -7653.727)
    ,@("GSM", "H76", "set", 35173)
    ,@("GSM", "J76", "set", 35173)
    ,@("GSM", "L76", "set", 35173)
    ,@("GSM", "N76", "set", -35803)
    ,@("GSM", "H79", "set", 35173)
    ,@("GSM", "J79", "set", 35173)
    ,@("GSM", "L79", "set", 35173)
    ,@("GSM", "N79", "set", -37357)
    ,@("GSM", "H93", "set", 22000)
    ,@("GSM", "I93", "set", 22000)
    ,@("GSM", "K93", "set", 22000)
    ,@("GSM", "M93", "set", -20128)
    ,@("GSM", "H99", "set", 12270)
    ,@("GSM", "I99", "set", 9693.333000000001)
    ,@("GSM", "J99", "set", 20000)
    ,@("GSM", "K99", "set", 9693.333000000001)
    ,@("GSM", "L99", "set", 20000)
    ,@("GSM", "M99", "set", -7447.333000000001)
    ,@("GSM", "N99", "set", -24492)
    ,@("GSM", "H122", "set", 3462.5417)
    ,@("GSM", "I122", "set", 3104.8096)
    ,@("GSM", "J122", "set", 5966.6665)
    ,@("GSM", "K122", "set", 9314.4288)
    ,@("GSM", "L122", "set", 17899.9995)
    ,@("GSM", "M122", "set", -6864.4288)
    ,@("GSM", "N122", "set", -22799.9995)
    ,@("GSM", "H132", "set", 58843.8)
    ,@("GSM", "I132", "set", 37015.5)
    ,@("GSM", "J132", "set", 146157)
    ,@("GSM", "K132", "set", 111046.5)
    ,@("GSM", "L132", "set", 438471)
    ,@("GSM", "M132", "set", -108516.5)
    ,@("GSM", "N132", "set", -443531)
    ,@("LTW", "H132", "set", 26184.46)
    ,@("LTW", "I132", "set", 11600.396)
    ,@("LTW", "J132", "set", 103480)
    ,@("LTW", "K132", "set", 34801.188)
    ,@("LTW", "L132", "set", 310440)
    ,@("LTW", "M132", "set", -32271.188)
    ,@("LTW", "N132", "set", -315500)
    ,@("LTW", "H141", "set", 59838.332)
    ,@("LTW", "J141", "set", 59838.332)
    ,@("LTW", "L141", "set", 59838.332)
    ,@("LTW", "N141", "set", -70198.33199999999)
    ,@("WVR", "H17", "set", 1250460)
    ,@("WVR", "I17", "set", 1250460)
    ,@("WVR", "K17", "set", 1250460)
    ,@("WVR", "M17", "set", -1250288)
    ,@("WVR", "H129", "set", 34985)
    ,@("WVR", "J129", "set", 34985)
    ,@("WVR", "L129", "set", 34985)
    ,@("WVR", "N129", "set", -44985)
    ,@("WVR", "H132", "set", 76399.67999999999)
    ,@("WVR", "I132", "set", 69696.44500000001)
    ,@("WVR", "J132", "set", 92599.164)
    ,@("WVR", "K132", "set", 209089.335)
    ,@("WVR", "L132", "set", 277797.492)
    ,@("WVR", "M132", "set", -206559.335)
    ,@("WVR", "N132", "set", -282857.492)
)

foreach ($edit in $edits) {
    $sheetName = $edit[0]
    $cellRef   = $edit[1]
    $action    = $edit[2]
    $newValue  = $edit[3]
    $ws = $wb.Worksheets.Item($sheetName)
    if ($action -eq "set") {
        $ws.Range($cellRef).Value = $newValue
    } else {
        $ws.Range($cellRef).Value = $null
    }
}

Write-Host "Applied $($edits.Count) cell updates."